# Update "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" worksheets, matching the upstream
# data refresh recorded in the commit "Update gh-pages to output
# generated at 456a3b4".
#
# NOTE: this runtime's Range/Cells ".Value" getter does not return the
# underlying scalar (it prints/compares as a reflection-like description),
# so all reads/compares below use ".Value2" instead, which behaves
# correctly for both get and set.

$wb = $excel.ActiveWorkbook

# column F (想去人数) updates on sheet "展览" - row -> (old, new)
$exhibitionUpdates = @(
    @{ Row = 2;  Old = 68;   New = 69   },
    @{ Row = 3;  Old = 525;  New = 527  },
    @{ Row = 4;  Old = 171;  New = 172  },
    @{ Row = 5;  Old = 222;  New = 228  },
    @{ Row = 6;  Old = 369;  New = 371  },
    @{ Row = 8;  Old = 2240; New = 2247 },
    @{ Row = 9;  Old = 378;  New = 380  },
    @{ Row = 10; Old = 5514; New = 5539 },
    @{ Row = 12; Old = 362;  New = 363  }
)

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($u in $exhibitionUpdates) {
    $cell = $wsExhibition.Cells.Item($u.Row, 6)  # column F
    if ($cell.Value2 -eq $u.Old) {
        $cell.Value2 = $u.New
    }
}

# column F (想去人数) updates on sheet "全部类型" - this sheet has two
# extra rows (rows 3 and 9/10) compared to "展览", so its row numbers
# for the same events are not a constant offset; list them explicitly.
$allTypesUpdates = @(
    @{ Row = 2;  Old = 68;   New = 69   },
    @{ Row = 4;  Old = 525;  New = 527  },
    @{ Row = 5;  Old = 171;  New = 172  },
    @{ Row = 6;  Old = 222;  New = 228  },
    @{ Row = 7;  Old = 369;  New = 371  },
    @{ Row = 11; Old = 2240; New = 2247 },
    @{ Row = 12; Old = 378;  New = 380  },
    @{ Row = 13; Old = 5514; New = 5539 },
    @{ Row = 15; Old = 362;  New = 363  }
)

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($u in $allTypesUpdates) {
    $cell = $wsAll.Cells.Item($u.Row, 6)  # column F
    if ($cell.Value2 -eq $u.Old) {
        $cell.Value2 = $u.New
    }
}
